# Updated documentation for MassWateR sampling method IDs and context
# (mirrors the "Add files via upload" commit to MassWateR_WQXMeta_Template.xlsx)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Meta")
$wsInstr = $wb.Worksheets.Item("Instructions")

# --- Meta tab -----------------------------------------------------------
# Column B ("Sampling Method Context") previously used the now-retired
# "MassBays" value for TP / TDP / E.coli rows; replace with "MassWateR".
$wsMeta.Range("B4").Value = "MassWateR"
$wsMeta.Range("B5").Value = "MassWateR"
$wsMeta.Range("B6").Value = "MassWateR"

# Re-apply alignment so the data rows share the same (unbordered) style
# that the rest of the workbook already uses, collapsing the bespoke
# bordered variants that existed only for this block.
$wsMeta.Range("A2:A6").HorizontalAlignment = -4131
$wsMeta.Range("B2:F6").HorizontalAlignment = -4108

# --- Instructions tab -----------------------------------------------------
# The row describing "Sampling Method Context" documents the new
# MassWateR-based context value (text unchanged; shared string index
# shifts once the unused "MassBays" entry is dropped).
$wsInstr.Range("B7").Value = "Enter the Context for the Sampling Method IDs that are used for sampling this parameter.  Not applicable for field measurements/observations.  If you are using the standard methods defined by MassWateR, enter the context ""MassWateR""."

# --- Selection / active-cell bookkeeping ---------------------------------
$wsInstr.Activate()
$wsInstr.Range("A7").Select()

$wsMeta.Activate()
$wsMeta.Range("A2").Select()
